$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4153.909
$ws.Range("J64").Value = 4224.25
$ws.Range("L64").Value = 4224.25
$ws.Range("N64").Value = -4720.25
$ws.Range("H67").Value = 4153.909
$ws.Range("J67").Value = 4224.25
$ws.Range("L67").Value = 4224.25
$ws.Range("N67").Value = -5940.25
$ws.Range("H74").Value = 10676.637
$ws.Range("I74").Value = 11382.777
$ws.Range("K74").Value = 11382.777
$ws.Range("M74").Value = -10446.777
$ws.Range("H77").Value = 10676.637
$ws.Range("I77").Value = 11382.777
$ws.Range("K77").Value = 56913.885
$ws.Range("M77").Value = -52233.885
$ws.Range("H92").Value = 585.3158
$ws.Range("I92").Value = 582.13336
$ws.Range("K92").Value = 582.13336
$ws.Range("M92").Value = 665.86664
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992
$ws.Range("H100").Value = 1366.6
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H116").Value = 4309.385
$ws.Range("I116").Value = 4736.4443
$ws.Range("J116").Value = 3348.5
$ws.Range("K116").Value = 4736.4443
$ws.Range("L116").Value = 3348.5
$ws.Range("M116").Value = -1294.4443
$ws.Range("N116").Value = -10232.5
$ws.Range("H125").Value = 1895895
$ws.Range("I125").Value = 3789545.2
$ws.Range("K125").Value = 34105906.8
$ws.Range("M125").Value = -34103446.8
$ws.Range("H138").Value = 5990.175
$ws.Range("I138").Value = 16596.6
$ws.Range("J138").Value = 2454.7
$ws.Range("K138").Value = 49789.8
$ws.Range("L138").Value = 7364.099999999999
$ws.Range("M138").Value = -44649.8
$ws.Range("N138").Value = -17644.1

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 776.7692
$ws.Range("I5").Value = 858
$ws.Range("J5").Value = 707.1429000000001
$ws.Range("K5").Value = 858
$ws.Range("L5").Value = 707.1429000000001
$ws.Range("M5").Value = -746
$ws.Range("N5").Value = -931.1429000000001
$ws.Range("H74").Value = 329269.47
$ws.Range("I74").Value = 1791.1177
$ws.Range("J74").Value = 1164339.2
$ws.Range("K74").Value = 1791.1177
$ws.Range("L74").Value = 1164339.2
$ws.Range("M74").Value = -917.1177
$ws.Range("N74").Value = -1166087.2
$ws.Range("H77").Value = 329269.47
$ws.Range("I77").Value = 1791.1177
$ws.Range("J77").Value = 1164339.2
$ws.Range("K77").Value = 8955.5885
$ws.Range("L77").Value = 5821696
$ws.Range("M77").Value = -4587.5885
$ws.Range("N77").Value = -5830432
$ws.Range("H97").Value = 6722.5
$ws.Range("I97").Value = 9291.385
$ws.Range("K97").Value = 9291.385
$ws.Range("M97").Value = -8795.385
$ws.Range("H102").Value = 1733.5
$ws.Range("I102").Value = 1733.5
$ws.Range("K102").Value = 1733.5
$ws.Range("M102").Value = -111.5
$ws.Range("H110").Value = 1944.1428
$ws.Range("I110").Value = 1944.1428
$ws.Range("K110").Value = 1944.1428
$ws.Range("M110").Value = 100.8571999999999
$ws.Range("H122").Value = 1650.25
$ws.Range("I122").Value = 1100.2858
$ws.Range("K122").Value = 3300.8574
$ws.Range("M122").Value = -850.8574000000003
$ws.Range("H132").Value = 2714.907
$ws.Range("I132").Value = 1675.6
$ws.Range("K132").Value = 5026.799999999999
$ws.Range("M132").Value = -2496.799999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 776.7692
$ws.Range("I4").Value = 858
$ws.Range("J4").Value = 707.1429000000001
$ws.Range("K4").Value = 858
$ws.Range("L4").Value = 707.1429000000001
$ws.Range("M4").Value = -743
$ws.Range("N4").Value = -937.1429000000001
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10566
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 3614.5
$ws.Range("I29").Value = 3152.6667
$ws.Range("K29").Value = 3152.6667
$ws.Range("M29").Value = -2863.6667
$ws.Range("H94").Value = 1554.3334
$ws.Range("I94").Value = 1633.7368
$ws.Range("K94").Value = 1633.7368
$ws.Range("M94").Value = -1182.7368

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1803.1666
$ws.Range("I16").Value = 877.5714
$ws.Range("K16").Value = 877.5714
$ws.Range("M16").Value = -590.5714
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H28").Value = 22646.25
$ws.Range("J28").Value = 22646.25
$ws.Range("L28").Value = 22646.25
$ws.Range("N28").Value = -23136.25
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H113").Value = 1803.1666
$ws.Range("I113").Value = 877.5714
$ws.Range("K113").Value = 877.5714
$ws.Range("M113").Value = 1292.4286

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 4067.75
$ws.Range("J50").Value = 3973.1428
$ws.Range("L50").Value = 11919.4284
$ws.Range("N50").Value = -12881.4284
$ws.Range("H53").Value = 4067.75
$ws.Range("J53").Value = 3973.1428
$ws.Range("L53").Value = 11919.4284
$ws.Range("N53").Value = -12881.4284
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H138").Value = 3993.3333
$ws.Range("I138").Value = 4175.9165
$ws.Range("K138").Value = 12527.7495
$ws.Range("M138").Value = -7387.749500000002
$ws.Range("H139").Value = 3135.4614
$ws.Range("I139").Value = 2341.3635
$ws.Range("K139").Value = 7024.0905
$ws.Range("M139").Value = -1884.0905
$ws.Range("H140").Value = 1609.2051
$ws.Range("I140").Value = 1233.0938
$ws.Range("K140").Value = 3699.2814
$ws.Range("M140").Value = 1480.7186
$ws.Range("H141").Value = 1961.5555
$ws.Range("I141").Value = 2081.875
$ws.Range("K141").Value = 6245.625
$ws.Range("M141").Value = -1065.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4104.25
$ws.Range("I122").Value = 3639.6667
$ws.Range("K122").Value = 10919.0001
$ws.Range("M122").Value = -8469.000100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1087.25
$ws.Range("I55").Value = 1336.9375
$ws.Range("K55").Value = 1336.9375
$ws.Range("M55").Value = -1163.9375
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H93").Value = 2804.25
$ws.Range("I93").Value = 1919.2142
$ws.Range("K93").Value = 1919.2142
$ws.Range("M93").Value = -671.2141999999999
$ws.Range("H95").Value = 49999.5
$ws.Range("J95").Value = 49999.5
$ws.Range("L95").Value = 49999.5
$ws.Range("N95").Value = -55491.5
$ws.Range("H105").Value = 14970.5
$ws.Range("J105").Value = 14970.5
$ws.Range("L105").Value = 14970.5
$ws.Range("N105").Value = -21958.5
$ws.Range("H132").Value = 4371
$ws.Range("I132").Value = 3150.625
$ws.Range("K132").Value = 9451.875
$ws.Range("M132").Value = -6921.875

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H96").Value = 35982.668
$ws.Range("I96").Value = 3999
$ws.Range("J96").Value = 67966.336
$ws.Range("K96").Value = 3999
$ws.Range("L96").Value = 67966.336
$ws.Range("M96").Value = -2626
$ws.Range("N96").Value = -70712.336
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 2283.8386
$ws.Range("I132").Value = 2036.5625
$ws.Range("J132").Value = 2547.6
$ws.Range("K132").Value = 6109.6875
$ws.Range("L132").Value = 7642.799999999999
$ws.Range("M132").Value = -3579.6875
$ws.Range("N132").Value = -12702.8
$ws.Range("H136").Value = 897.8461
$ws.Range("I136").Value = 736.125
$ws.Range("J136").Value = 1156.6
$ws.Range("K136").Value = 2208.375
$ws.Range("L136").Value = 3469.8
$ws.Range("M136").Value = 341.625
$ws.Range("N136").Value = -8569.799999999999
